$d = $word.ActiveDocument

$replacements = @(
    ,@("N = 90,190", "N = 34,218")
    ,@("96.4 (89.7, 102.7)", "96.4 (89.7, 102.8)")
    ,@("89,875 (100)", "34,100 (100)")
    ,@("89,545 (99)", "33,974 (99)")
    ,@("463.8 (319.4, 641.1)", "466.8 (323.3, 641.4)")
    ,@("234.0 (114.5, 404.0)", "236.1 (117.0, 410.0)")
    ,@("585.7 (372.0, 852.8)", "589.5 (375.4, 856.6)")
    ,@("2,016 (2.2)", "773 (2.3)")
    ,@("1,282 (1.4)", "497 (1.5)")
    ,@("63.2 (56.1, 68.4)", "63.6 (56.6, 68.6)")
    ,@("2,677 (3.0)", "1,025 (3.0)")
    ,@("87,513 (97)", "33,193 (97)")
    ,@("51,685 (57)", "19,648 (57)")
    ,@("38,505 (43)", "14,570 (43)")
    ,@("7,103 (7.9)", "2,700 (7.9)")
    ,@("22,101 (25)", "8,281 (24)")
    ,@("21,308 (24)", "8,044 (24)")
    ,@("39,678 (44)", "15,193 (44)")
    ,@("-2.5 (-3.8, -0.2)", "-2.5 (-3.9, -0.3)")
    ,@("11,568 (13)", "4,418 (13)")
    ,@("19,649 (22)", "7,562 (22)")
    ,@("23,651 (26)", "8,952 (26)")
    ,@("20,891 (23)", "7,825 (23)")
    ,@("6,072 (6.7)", "2,248 (6.6)")
    ,@("8,359 (9.3)", "3,213 (9.4)")
    ,@("512 (0.6)", "197 (0.6)")
    ,@("35,505 (39)", "13,588 (40)")
    ,@("37,129 (41)", "14,175 (41)")
    ,@("17,044 (19)", "6,258 (18)")
    ,@("87,335 (97)", "33,086 (97)")
    ,@("18,184 (21)", "6,824 (21)")
    ,@("69,482 (79)", "26,451 (79)")
    ,@("52,073 (58)", "19,710 (58)")
    ,@("32,043 (36)", "12,213 (36)")
    ,@("6,074 (6.7)", "2,295 (6.7)")
    ,@("4,954 (5.5)", "1,879 (5.5)")
    ,@("18,230 (20)", "6,826 (20)")
    ,@("22,684 (25)", "8,456 (25)")
    ,@("23,623 (26)", "8,987 (26)")
    ,@("20,699 (23)", "8,070 (24)")
    ,@("64,233 (71)", "24,361 (71)")
    ,@("22,761 (25)", "8,664 (25)")
    ,@("3,196 (3.5)", "1,193 (3.5)")
    ,@("16,282 (18)", "6,177 (18)")
    ,@("30,576 (34)", "11,667 (34)")
    ,@("43,332 (48)", "16,374 (48)")
    ,@("74,487 (83)", "28,275 (83)")
    ,@("14,982 (17)", "5,675 (17)")
    ,@("721 (0.8)", "268 (0.8)")
    ,@("76,495 (85)", "29,083 (85)")
    ,@("13,271 (15)", "4,979 (15)")
    ,@("424 (0.5)", "156 (0.5)")
    ,@("29,441 (33)", "11,115 (32)")
    ,@("54,016 (60)", "20,568 (60)")
    ,@("6,733 (7.5)", "2,535 (7.4)")
    ,@("15,354 (17)", "5,761 (17)")
    ,@("66,435 (74)", "25,213 (74)")
    ,@("8,401 (9.3)", "3,244 (9.5)")
    ,@("0.1 (0.0, 0.6)", "0.1 (0.0, 0.7)")
    ,@("29.8 (16.8, 50.3)", "29.9 (16.8, 50.6)")
    ,@("19.4 (6.0, 41.1)", "19.6 (6.0, 41.5)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"